$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'46.516.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  +1.74%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Formula = "'2.621.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +8.51%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Formula = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Formula = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Formula = "'314.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  +4.48%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Formula = "'102.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  +4.45%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Formula = "'0.603"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  +6.61%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Formula = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'  +0.03%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Formula = "'0.592"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'  +14.83%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Formula = "'39.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  +11.19%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Formula = "'54.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  +1.05%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Formula = "'0.0844"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  +6.11%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Formula = "'  +16.04%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Formula = "'3.020.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  +9.00%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Formula = "'  +1.74%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Formula = "'2.624.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  +8.72%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Formula = "'0.918"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  +8.05%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Formula = "'15.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  +5.79%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Formula = "'46.768.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  +2.37%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Formula = "'13.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  +0.98%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Formula = "'  +7.87%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Formula = "'  +8.46%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Formula = "'71.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  +5.36%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Formula = "'255.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  +4.79%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Formula = "'3.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  +10.65%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Formula = "'  +14.35%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Formula = "'28.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "'  +32.08%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Formula = "'  +0.02%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Formula = "'InjectiveProtocol"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Formula = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Formula = "'41.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "'  +6.62%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Formula = "'Cosmos"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Formula = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Formula = "'10.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  +8.89%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Formula = "'2.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  +2.35%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Formula = "'  +11.38%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Formula = "'3.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  -2.68%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Formula = "'  +13.64%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Formula = "'2.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "'  +4.11%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Formula = "'154.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  +3.68%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Formula = "'0.0840"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  +7.92%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Formula = "'  +4.66%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Formula = "'  +6.00%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Formula = "'17.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  +11.82%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Formula = "'4.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  +10.16%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Formula = "'3.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  +11.01%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Formula = "'0.0330"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  +9.39%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Formula = "'21.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'  +37.56%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Formula = "'2.040.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  +4.35%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Formula = "'  -0.06%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Formula = "'91.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  -0.13%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Formula = "'113.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  +10.07%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Formula = "'  +3.24%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Formula = "'9.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  +6.48%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Formula = "'78.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +12.89%  "
$ws.Range("E51").Style = "Normal"
